$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert this week's new record at row 229 (pushing the
# existing history for "Feria Lagunitas de Puerto Montt - Pepino ensalada"
# down by one row) and populate it with the latest observation.
$ws.Rows("229").Insert()

$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44754
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = 100112043
$ws.Range("G229").Value = "Pepino ensalada"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 400
$ws.Range("K229").Value = 23000
$ws.Range("L229").Value = 23000
$ws.Range("M229").Value = 23000
$ws.Range("N229").Value = "$/caja 60 unidades"
$ws.Range("O229").Value = "Región de Arica y Parinacota"
$ws.Range("P229").Value = 383
$ws.Range("Q229").Value = 60
$ws.Range("R229").Value = "Hortaliza"
